$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 20: Intake's "List of possible next queues" now also allows Suspend
$ws.Range("D20").Value = "Fulfill,Suspend"

# Row 24: Release's "Default next queue" becomes "Executive Approval"
$ws.Range("E24").Value = """Executive Approval"""

# New Row 25: Suspend Next Queues
$ws.Range("B25").Value = "Suspend Next Queues"
$ws.Range("C25").Value = "Suspend"
$ws.Range("D25").Value = "Fulfill"
$ws.Range("E25").Value = """Fulfill"""
$ws.Range("F25").Value = "null"
